# generated copilot script test
# Rework the "MASTER WORKSHEET" sheet: remove the old junk header columns,
# promote the device-list column (formerly column B) into column A, and
# rename the header cell to "Device Names".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MASTER WORKSHEET")

# Drop the old column A ("Worksheet = master worksheet" / "Junk A2 Text");
# this shifts the device-name column (old column B) left into column A,
# carrying over its column width/formatting.
$ws.Columns.Item(1).Delete()

# Drop the old header row 1 (now holding "Column B Header Text" /
# "Device Name 1"), shifting the CAM/GSS/KEY/INP device rows up.
$ws.Rows.Item(1).Delete()

# Rename the remaining header cell.
$ws.Range("A1").Value = "Device Names"

# Restore cursor position recorded for this sheet.
[void]$ws.Range("G14").Select()
